$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "mm"
$ws.Range("A2").Value = "class"
$ws.Range("B2").Value = "A"

$ws.Range("B2").Select()
